
# ---------------------------------------------------------------------------
# Generate Report for Handoff
#
# Adds a new row (row 3) to each of the three worksheets ("Overview",
# "zh-cn", "de-de") describing the file
#   598141e5-97f8-4537-aa32-399c8d758868ooo....md
# which has now reached the "Ready for handoff" state. Mirrors the existing
# row 2 (the 2b856dbb... file) on every sheet: duplicates formatting via a
# row copy/insert, then overwrites the cells that differ, re-adds the
# hyperlink, resizes each table (ListObject) to include the new row, and
# widens the two "latest handoff/handback file" columns slightly to fit the
# new, longer status text ("Ready for handoff").
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$newMdName    = '598141e5-97f8-4537-aa32-399c8d758868ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$newMdDisplay = 'e2e\598141e5-97f8-4537-aa32-399c8d758868ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$readyStatus  = 'Ready for handoff'
$overviewDate = '2016-09-05 20:33:53'
$xlfZhCnNew   = '598141e5-97f8-4537-aa32-399c8d758868oooooooooooooooooooooooooooooooooooooooo.cfac9e91bdc0c77dfdb061ca4eb841014e2c581b.zh-cn.xlf'
$zhCnHbDate   = '2016-09-05 20:33:48'
$xlfDeDeNew   = '598141e5-97f8-4537-aa32-399c8d758868oooooooooooooooooooooooooooooooooooooooo.cfac9e91bdc0c77dfdb061ca4eb841014e2c581b.de-de.xlf'
$newMdUrl     = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bef53f1bf32221591b1ce01ba429eae98c59bed5/e2e/598141e5-97f8-4537-aa32-399c8d758868ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'

# ===========================================================================
# Sheet "Overview" : File Name | Path And Name | Extension | Publish URL |
#                    zh-cn | de-de | Latest HO Xliff Generate Date
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows(2).Copy()
$wsOverview.Rows(3).Insert()

$wsOverview.Range("A3").Value = $newMdName
$wsOverview.Range("B3").Value = $newMdDisplay
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = $overviewDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newMdUrl, "", "", $newMdDisplay)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

$wsOverview.Columns(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns(6).ColumnWidth = 16.333333333333336

# ===========================================================================
# Sheet "zh-cn" : Source File Name | File Extension | Status | Source Path |
#   Priority | Content Duplicate | Latest Handoff File | Latest Handoff
#   Datetime | Latest Target File | Latest Handback File | Latest Handback
#   DateTime | Reference Tokens | To be localized | Dependency From |
#   Has metadata | Error Detail
# ===========================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Rows(2).Copy()
$wsZhCn.Rows(3).Insert()

$wsZhCn.Range("A3").Value = $newMdName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $xlfZhCnNew
$wsZhCn.Range("H3").Value = $zhCnHbDate
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newMdUrl, "", "", $newMdDisplay)

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

$wsZhCn.Columns(3).ColumnWidth = 16.333333333333336

# ===========================================================================
# Sheet "de-de" : same columns as "zh-cn"
# ===========================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Rows(2).Copy()
$wsDeDe.Rows(3).Insert()

$wsDeDe.Range("A3").Value = $newMdName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $xlfDeDeNew
$wsDeDe.Range("H3").Value = $overviewDate
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newMdUrl, "", "", $newMdDisplay)

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

$wsDeDe.Columns(3).ColumnWidth = 16.333333333333336

Write-Host "Done: added handoff row to Overview, zh-cn, de-de sheets."
